# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price ("D") and Volume(1h) ("E") columns are plain text in this sheet (not
# numbers), so every "D" write forces Text format first and clears the
# resulting style afterwards -- this keeps the numeric-looking price strings
# (e.g. "619.50", "1.00", "0.0000248") stored verbatim as text instead of
# being auto-coerced into Excel numbers, while leaving no extra cell styling
# behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.554.96'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.211.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +4.38%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("E5").Value = '  +0.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.50'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("E7").Value = '  +5.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.372'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.36%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.216.78'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.734'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.44%  '

$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("B13").Value = 'WrappedBTC'
$ws.Range("C13").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '101.382.97'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +13.32%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.15%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.42'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.57'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.815.28'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.221.34'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.10%  '

$ws.Range("E19").Value = '  -2.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.22'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +10.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +11.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '451.81'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000204'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.25'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.02'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +8.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.330.95'
$ws.Range("D28").ClearFormats()

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.144'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +59.05%  '

$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.233'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +18.43%  '

$ws.Range("E32").Value = '  +7.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.40'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.171'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +13.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.18'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.68'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +7.47%  '

$ws.Range("E38").Value = '  +5.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '511.16'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.06%  '

$ws.Range("E40").Value = '  +7.65%  '

$ws.Range("E41").Value = '  +13.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.80'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -9.30%  '

$ws.Range("E43").Value = '  -4.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.08'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.740'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.94'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.51'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("E49").Value = '  +7.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.53'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.96'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.98%  '
